# تعديل تلقائي في شيت Card7 by admin at 2025-12-06 18:32:34
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card7")

for ($r = 3; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    # Force text storage (matches existing "card" column cells, which are
    # numeric-looking strings) instead of letting Excel auto-convert "7"
    # into a Number, then restore the default "Normal" style so no stray
    # number-format/style is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = "7"
    $cell.Style = "Normal"
}
